# Merge the split "<id>...</id>" runs back into a single run for both
# the p128r_1 and p128v_1 identifier paragraphs.
#
# Before: three runs -> "<id>", "p128r_1", "</id>"
# After:  one run    -> "<id>p128r_1</id>"
# (and similarly for p128v_1)
#
# Using Find/Replace on the plain text that spans the three runs causes
# Word to collapse them into a single run that carries the formatting of
# the first (leading) run, which is exactly the formatting kept in the
# target document.

$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p128r_1</id>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p128r_1</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p128v_1</id>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p128v_1</id>", 2) | Out-Null
